{"js": "// PerformanceLevelRiskPageTemplate edit:\n//   1. Remove the old \"_GoBack\" bookmark (it used to sit right before the\n//      \"PL\" run in the \"From initial risk ... PL\" heading line).\n//   2. On the final line of the document (\"(PLr)Performance level: <PerformanceLevel>\")\n//      insert the text \"Required \" right after \"(PLr)\" and place a fresh\n//      \"_GoBack\" bookmark right after that inserted text (i.e. immediately\n//      before \"Performance level: \").\n\nconst body = context.document.body;\n\n// --- Step 1: remove the pre-existing _GoBack bookmark -------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Step 2: insert \"Required \" after \"(PLr)\" and re-create the bookmark\nconst results = body.search(\"(PLr)\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n    const plrRange = results.items[0];\n\n    // Insert the new run right after \"(PLr)\"; the inherited formatting\n    // (bold, Gotham Light, en-GB) matches the run that used to hold\n    // \"Performance level: \".\n    const insertedRange = plrRange.insertText(\"Required \", Word.InsertLocation.after);\n    await context.sync();\n\n    // Collapse to the point right after the inserted text, then drop the\n    // \"_GoBack\" bookmark there (i.e. just before \"Performance level: \").\n    const afterInserted = insertedRange.getRange(Word.RangeLocation.after);\n    afterInserted.insertBookmark(\"_GoBack\");\n    await context.sync();\n}\n", "ps1": "# PerformanceLevelRiskPageTemplate edit:\n#   1. Remove the old \"_GoBack\" bookmark (it used to sit right before the\n#      \"PL\" run in the \"From initial risk ... PL\" heading line).\n#   2. On the final line of the document (\"(PLr)Performance level: <PerformanceLevel>\")\n#      insert the text \"Required \" right after \"(PLr)\" and place a fresh\n#      \"_GoBack\" bookmark right after that inserted text (i.e. immediately\n#      before \"Performance level: \").\n\n$d = $word.ActiveDocument\n\n# --- Step 1: remove the pre-existing _GoBack bookmark -----------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- Step 2: insert \"Required \" after \"(PLr)\" and re-create the bookmark\n$rng = $d.Content\n$found = $rng.Find.Execute(\"(PLr)\")\nif ($found) {\n    $rng.Collapse(0)               # wdCollapseEnd - collapse to an insertion point right after \"(PLr)\"\n    $rng.InsertAfter(\"Required \")  # type the new text\n    $rng.Collapse(0)               # collapse again so the bookmark doesn't wrap \"Required \"\n    $d.Bookmarks.Add(\"_GoBack\", $rng)\n}\n"}
